$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 54233.49
$ws.Range("I33").Value = 30470.818
$ws.Range("J33").Value = 250275.5
$ws.Range("K33").Value = 30470.818
$ws.Range("L33").Value = 250275.5
$ws.Range("M33").Value = -30241.818
$ws.Range("N33").Value = -250733.5
$ws.Range("H115").Value = 5097.5
$ws.Range("I115").Value = 2117
$ws.Range("K115").Value = 6351
$ws.Range("M115").Value = -4784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1575.8667
$ws.Range("I45").Value = 1577
$ws.Range("J45").Value = 1571.3334
$ws.Range("K45").Value = 1577
$ws.Range("L45").Value = 1571.3334
$ws.Range("M45").Value = -1200
$ws.Range("N45").Value = -2325.3334
$ws.Range("H61").Value = 1103.1333
$ws.Range("I61").Value = 1071.5358
$ws.Range("J61").Value = 1545.5
$ws.Range("K61").Value = 1071.5358
$ws.Range("L61").Value = 1545.5
$ws.Range("M61").Value = -859.5358000000001
$ws.Range("N61").Value = -1969.5
$ws.Range("H74").Value = 1139.1296
$ws.Range("I74").Value = 725.4706
$ws.Range("J74").Value = 8171.3335
$ws.Range("K74").Value = 725.4706
$ws.Range("L74").Value = 8171.3335
$ws.Range("M74").Value = 148.5294
$ws.Range("N74").Value = -9919.333500000001
$ws.Range("H77").Value = 1139.1296
$ws.Range("I77").Value = 725.4706
$ws.Range("J77").Value = 8171.3335
$ws.Range("K77").Value = 3627.353
$ws.Range("L77").Value = 40856.6675
$ws.Range("M77").Value = 740.6469999999999
$ws.Range("N77").Value = -49592.6675
$ws.Range("H101").Value = 40515.855
$ws.Range("J101").Value = 40515.855
$ws.Range("L101").Value = 40515.855
$ws.Range("N101").Value = -47005.855
$ws.Range("H122").Value = 1673
$ws.Range("I122").Value = 1196.125
$ws.Range("J122").Value = 2626.75
$ws.Range("K122").Value = 3588.375
$ws.Range("L122").Value = 7880.25
$ws.Range("M122").Value = -1138.375
$ws.Range("N122").Value = -12780.25
$ws.Range("H132").Value = 7851.2104
$ws.Range("I132").Value = 9144.933999999999
$ws.Range("J132").Value = 2999.75
$ws.Range("K132").Value = 27434.802
$ws.Range("L132").Value = 8999.25
$ws.Range("M132").Value = -24904.802
$ws.Range("N132").Value = -14059.25
$ws.Range("H136").Value = 1103.1333
$ws.Range("I136").Value = 1071.5358
$ws.Range("J136").Value = 1545.5
$ws.Range("K136").Value = 3214.6074
$ws.Range("L136").Value = 4636.5
$ws.Range("M136").Value = -664.6074000000003
$ws.Range("N136").Value = -9736.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4215.3887
$ws.Range("I134").Value = 4483.032
$ws.Range("J134").Value = 2556
$ws.Range("K134").Value = 13449.096
$ws.Range("L134").Value = 7668
$ws.Range("M134").Value = -10914.096
$ws.Range("N134").Value = -12738

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1127.8
$ws.Range("I58").Value = 1140.725
$ws.Range("J58").Value = 1076.1
$ws.Range("K58").Value = 1140.725
$ws.Range("L58").Value = 1076.1
$ws.Range("M58").Value = -937.7249999999999
$ws.Range("N58").Value = -1482.1
$ws.Range("H99").Value = 1468.3125
$ws.Range("I99").Value = 1432.8182
$ws.Range("J99").Value = 1546.4
$ws.Range("K99").Value = 1432.8182
$ws.Range("L99").Value = 1546.4
$ws.Range("M99").Value = 65.18180000000007
$ws.Range("N99").Value = -4542.4
$ws.Range("H126").Value = 1468.3125
$ws.Range("I126").Value = 1432.8182
$ws.Range("J126").Value = 1546.4
$ws.Range("K126").Value = 4298.4546
$ws.Range("L126").Value = 4639.200000000001
$ws.Range("M126").Value = -1828.4546
$ws.Range("N126").Value = -9579.200000000001
$ws.Range("H132").Value = 10028.792
$ws.Range("J132").Value = 17144.572
$ws.Range("L132").Value = 51433.716
$ws.Range("N132").Value = -56493.716
$ws.Range("H134").Value = 6562
$ws.Range("I134").Value = 1694.5
$ws.Range("J134").Value = 26032
$ws.Range("K134").Value = 5083.5
$ws.Range("L134").Value = 78096
$ws.Range("M134").Value = -2548.5
$ws.Range("N134").Value = -83166
$ws.Range("H136").Value = 1127.8
$ws.Range("I136").Value = 1140.725
$ws.Range("J136").Value = 1076.1
$ws.Range("K136").Value = 3422.175
$ws.Range("L136").Value = 3228.3
$ws.Range("M136").Value = -872.1749999999997
$ws.Range("N136").Value = -8328.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2452
$ws.Range("I136").Value = 1986.6666
$ws.Range("J136").Value = 3150
$ws.Range("K136").Value = 5959.9998
$ws.Range("L136").Value = 9450
$ws.Range("M136").Value = -859.9997999999996
$ws.Range("N136").Value = -19650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5146
$ws.Range("I70").Value = 4989.8237
$ws.Range("J70").Value = 5323
$ws.Range("K70").Value = 4989.8237
$ws.Range("L70").Value = 5323
$ws.Range("M70").Value = -4719.8237
$ws.Range("N70").Value = -5863
$ws.Range("H73").Value = 5146
$ws.Range("I73").Value = 4989.8237
$ws.Range("J73").Value = 5323
$ws.Range("K73").Value = 4989.8237
$ws.Range("L73").Value = 5323
$ws.Range("M73").Value = -4053.8237
$ws.Range("N73").Value = -7195
$ws.Range("H80").Value = 2261.818
$ws.Range("I80").Value = 1970
$ws.Range("J80").Value = 2428.5715
$ws.Range("K80").Value = 1970
$ws.Range("L80").Value = 2428.5715
$ws.Range("M80").Value = -972
$ws.Range("N80").Value = -4424.5715
$ws.Range("H83").Value = 2261.818
$ws.Range("I83").Value = 1970
$ws.Range("J83").Value = 2428.5715
$ws.Range("K83").Value = 9850
$ws.Range("L83").Value = 12142.8575
$ws.Range("M83").Value = -4858
$ws.Range("N83").Value = -22126.8575
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H132").Value = 4449.1055
$ws.Range("I132").Value = 4396
$ws.Range("J132").Value = 4799.6
$ws.Range("K132").Value = 13188
$ws.Range("L132").Value = 14398.8
$ws.Range("M132").Value = -10658
$ws.Range("N132").Value = -19458.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1721.1111
$ws.Range("I7").Value = 1034
$ws.Range("J7").Value = 1985.3846
$ws.Range("K7").Value = 1034
$ws.Range("L7").Value = 1985.3846
$ws.Range("M7").Value = -922
$ws.Range("N7").Value = -2209.3846
$ws.Range("H75").Value = 2000
$ws.Range("I75").Value = 2000
$ws.Range("K75").Value = 2000
$ws.Range("M75").Value = -1064
$ws.Range("H78").Value = 2000
$ws.Range("I78").Value = 2000
$ws.Range("K78").Value = 6000
$ws.Range("M78").Value = -1320
$ws.Range("H104").Value = 9412.714
$ws.Range("J104").Value = 9412.714
$ws.Range("L104").Value = 9412.714
$ws.Range("N104").Value = -16400.714
$ws.Range("H126").Value = 1721.1111
$ws.Range("I126").Value = 1034
$ws.Range("J126").Value = 1985.3846
$ws.Range("K126").Value = 3102
$ws.Range("L126").Value = 5956.1538
$ws.Range("M126").Value = -632
$ws.Range("N126").Value = -10896.1538
$ws.Range("H136").Value = 2350.8076
$ws.Range("J136").Value = 3126.7307
$ws.Range("L136").Value = 9380.1921
$ws.Range("N136").Value = -14480.1921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2288.3467
$ws.Range("I132").Value = 2391.9692
$ws.Range("J132").Value = 1614.8
$ws.Range("K132").Value = 7175.9076
$ws.Range("L132").Value = 4844.4
$ws.Range("M132").Value = -4645.9076
$ws.Range("N132").Value = -9904.4
$ws.Range("H136").Value = 1407.625
$ws.Range("I136").Value = 1444.0714
$ws.Range("J136").Value = 1152.5
$ws.Range("K136").Value = 4332.2142
$ws.Range("L136").Value = 3457.5
$ws.Range("M136").Value = -1782.2142
$ws.Range("N136").Value = -8557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N110").ClearContents()
